$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")
$ws.Select()

$ws.Range("A3").Value = "juanosorio199@gmail.com"
